$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.092207074165344
$ws.Range("B1").Value = 1.01624608039856
$ws.Range("C1").Value = 0.8136826753616333
$ws.Range("D1").Value = 0.8097929954528809
$ws.Range("E1").Value = 0.8984602093696594
